# ---------------------------------------------------------------------------
# Apply the "Wait changes -- implicit" commit:
#   * Capabilities sheet: refresh the pCloudyAndroid credentials/app in row 6
#     (new uploader email, session token, apk name) and drop the now-removed
#     hyperlink on M6 (M7's hyperlink is kept).
#   * DeviceList sheet: grow the single-device table (column B) into a
#     10-device table (columns B:K) with the new Samsung device fleet.
#   * Test Data sheet: add a header row listing the next 10 Samsung devices.
#   * Misc view bookkeeping: selections / active sheet to match the saved file.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Capabilities sheet (sheet1)
# ---------------------------------------------------------------------------
$wsCap = $wb.Worksheets.Item("Capabilities")

# Row 6 (pCloudyAndroid row): new pCloudy username / api key / application.
$wsCap.Range("M6").Value = "sakshi.juneja@crestechsoftware.com"
$wsCap.Range("N6").Value = "t68k6kw68ywjv2y9zwfr9r3t"
$wsCap.Range("O6").Value = "dbsandroidUAT_10Nov.apk"

# The M6 mailto hyperlink is gone in the new version; M7 keeps its link.
$m6Addr = '$M$6'
foreach ($link in $wsCap.Hyperlinks) {
    if ($link.Range.Address() -eq $m6Addr) {
        $link.Delete()
    }
}

# ---------------------------------------------------------------------------
# 2. DeviceList sheet (sheet2) -- expand from one device column to ten.
# ---------------------------------------------------------------------------
$wsDev = $wb.Worksheets.Item("DeviceList")

$devCols = @("B", "C", "D", "E", "F", "G", "H", "I", "J", "K")

$deviceNames = @(
    "SAMSUNG_GalaxyJ6_android_10.0.0_482da",
    "SAMSUNG_GalaxyJ7Pro_android_9.0.0_a715a",
    "SAMSUNG_GalaxyJ7Pro_android_8.1.0_42e4e",
    "SAMSUNG_GalaxyJ8_android_10.0.0_882d2",
    "SAMSUNG_GalaxyM01_android_11.0.0_7425f",
    "SAMSUNG_GalaxyM02_android_11.0.0_51323",
    "SAMSUNG_GalaxyM10_android_10.0.0_a58e4",
    "SAMSUNG_GalaxyM12_android_11.0.0_df6a7",
    "SAMSUNG_GalaxyM20_android_8.1.0_258bd",
    "SAMSUNG_GalaxyM40_android_11.0.0_efd9b"
)

$versions = @("11.0.0","11.0.0","11.0.0","11.0.0","11.0.0","11.0.0","11.0.0","10.0.0","9.0.0","10.0.0")

$userNames = @(
    "S2021218GUID","S2021219EUID","S2021220IUID","S2325474CUID","S2325475AUID",
    "S2325476ZUID","S2325477HUID","S2325478FUID","S2325479DUID","S2325480HUID"
)

for ($i = 0; $i -lt $devCols.Length; $i++) {
    $col = $devCols[$i]
    $wsDev.Range("$col`1").Value = $deviceNames[$i]
    $wsDev.Range("$col`2").Value = $versions[$i]
    $wsDev.Range("$col`3").Value = "pCloudyAndroid"
    $wsDev.Range("$col`4").Value = $userNames[$i]
    $wsDev.Range("$col`5").Value = "121212"
    $wsDev.Range("$col`6").Value = "Xiaomi"
    $wsDev.Range("$col`7").Value = "10.0.0"
    $wsDev.Range("$col`8").Value = "12.0.0"
    $wsDev.Range("$col`9").Value = ($i + 1)
    $wsDev.Range("$col`10").Value = "DBS"
}

# Column widths for the newly-populated columns.
$wsDev.Columns.Item(2).ColumnWidth = 15.42578125
$wsDev.Columns.Item(3).ColumnWidth = 16.7109375
$wsDev.Columns.Item(7).ColumnWidth = 17.140625
$wsDev.Columns.Item(8).ColumnWidth = 16.85546875
$wsDev.Columns.Item(9).ColumnWidth = 15.28515625
$wsDev.Columns.Item(10).ColumnWidth = 17
$wsDev.Columns.Item(11).ColumnWidth = 14.7109375

$wsDev.Rows.Item(1).RowHeight = 45
$wsDev.Range("A1:K1").WrapText = $true
$wsDev.Range("A1:K1").Borders.LineStyle = 1

# ---------------------------------------------------------------------------
# 3. Test Data sheet (sheet3) -- header row with the next batch of devices.
# ---------------------------------------------------------------------------
$wsTD = $wb.Worksheets.Item("Test Data")

$tdCols = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J")
$tdDeviceNames = @(
    "SAMSUNG_GalaxyA10s_android_11.0.0_09401",
    "SAMSUNG_GalaxyA12_android_11.0.0_334bc",
    "SAMSUNG_GalaxyA21s_android_11.0.0_b13a4",
    "SAMSUNG_GalaxyA31_android_11.0.0_32c0a",
    "SAMSUNG_GalaxyA50_android_11.0.0_310bf",
    "SAMSUNG_GalaxyA51_android_11.0.0_d52ba",
    "SAMSUNG_GalaxyA71_android_11.0.0_fe4a3",
    "SAMSUNG_GalaxyA9_android_10.0.0_6eed1",
    "SAMSUNG_GalaxyFold_android_9.0.0_d69de",
    "SAMSUNG_GalaxyJ4_android_10.0.0_889f0"
)

for ($i = 0; $i -lt $tdCols.Length; $i++) {
    $col = $tdCols[$i]
    $wsTD.Range("$col`1").Value = $tdDeviceNames[$i]
}

$wsTD.Rows.Item(1).RowHeight = 90
$wsTD.Range("A1:J1").WrapText = $true
$wsTD.Range("A1:J1").Borders.LineStyle = 1

# ---------------------------------------------------------------------------
# 4. View bookkeeping -- selections on each sheet, re-activating DeviceList
#    last so it stays the saved "active tab" (matches the original file).
# ---------------------------------------------------------------------------
$wsCap.Range("B13").Select()
$wsTD.Range("F9").Select()
$wsDev.Range("D11").Select()
